# Generate Report for Handoff
# Updates the "Latest Handoff Date(time)" values to reflect a fresh report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1. Overview sheet: "Latest Handoff Date" column (D), rows 2-10
#    refreshed from 2016-03-23 05:16:25 -> 2016-03-23 05:17:12
$wsOverview.Range("D2:D10").Value = "2016-03-23 05:17:12"

# 2. de-de sheet: "Latest Handoff Datetime" column (E), rows 2-10
#    same refreshed timestamp as the Overview sheet (shares the same source value)
$wsDeDe.Range("E2:E10").Value = "2016-03-23 05:17:12"

# 3. zh-cn sheet: "Latest Handoff Datetime" column (E), rows 4-10
#    refreshed from 2016-03-23 05:16:13 -> 2016-03-23 05:17:02
#    (rows 2-3 already reflect a later/independent handoff and are left untouched)
$wsZhCn.Range("E4:E10").Value = "2016-03-23 05:17:02"
